{"js": "const pairs = [\n  [\"53-23=30\", \"30-4=26\"],\n  [\"76-35=41\", \"6+71=77\"],\n  [\"35+55=90\", \"81-56=25\"],\n  [\"20+38=58\", \"93-27=66\"],\n  [\"40-16=24\", \"77-55=22\"],\n  [\"39+8=47\", \"71-3=68\"],\n  [\"73-33=40\", \"43-14=29\"],\n  [\"2+56=58\", \"87+5=92\"],\n  [\"4+92=96\", \"99-78=21\"],\n  [\"30-29=1\", \"16+60=76\"],\n  [\"95-56=39\", \"56-50=6\"],\n  [\"76-21=55\", \"76-60=16\"],\n  [\"74-12=62\", \"84-14=70\"],\n  [\"57-10=47\", \"23+62=85\"],\n  [\"50-15=35\", \"12+74=86\"],\n  [\"22+29=51\", \"18+19=37\"],\n  [\"95-12=83\", \"37+55=92\"],\n  [\"70+11=81\", \"85-30=55\"],\n  [\"24+34=58\", \"97-62=35\"],\n  [\"41+53=94\", \"92-55=37\"],\n  [\"73-1=72\", \"4+65=69\"],\n  [\"63+33=96\", \"61-20=41\"],\n  [\"43-13=30\", \"84-3=81\"],\n  [\"13+31=44\", \"50-45=5\"],\n  [\"57-22=35\", \"26+48=74\"],\n  [\"67-30=37\", \"39-31=8\"],\n  [\"1+42=43\", \"88-57=31\"],\n  [\"12+72=84\", \"17+56=73\"],\n  [\"87-82=5\", \"16+71=87\"],\n  [\"75+11=86\", \"11+72=83\"],\n  [\"7+53=60\", \"93-9=84\"],\n  [\"60+35=95\", \"60-38=22\"],\n  [\"64-51=13\", \"72-22=50\"],\n  [\"66-44=22\", \"5+6=11\"],\n  [\"46-32=14\", \"88-21=67\"],\n  [\"20-3=17\", \"21+53=74\"],\n  [\"12+86=98\", \"43+11=54\"],\n  [\"4+25=29\", \"7+11=18\"],\n  [\"38+23=61\", \"11+45=56\"],\n  [\"33-16=17\", \"28-9=19\"],\n  [\"69+1=70\", \"39-5=34\"],\n  [\"28+25=53\", \"28+29=57\"],\n  [\"58-55=3\", \"41-7=34\"],\n  [\"22+11=33\", \"58-54=4\"],\n  [\"21+70=91\", \"69-32=37\"],\n  [\"31-25=6\", \"30+11=41\"],\n  [\"66-18=48\", \"43+20=63\"],\n  [\"53+28=81\", \"31+21=52\"],\n  [\"5+58=63\", \"64-38=26\"],\n  [\"18+38=56\", \"51-27=24\"],\n  [\"60+26=86\", \"42+0=42\"],\n  [\"19+45=64\", \"63-20=43\"],\n  [\"64+14=78\", \"46-6=40\"],\n  [\"26+30=56\", \"7+26=33\"],\n  [\"47-30=17\", \"53-33=20\"],\n  [\"68-38=30\", \"10+22=32\"],\n  [\"36-20=16\", \"12+36=48\"],\n  [\"29+59=88\", \"28+0=28\"],\n  [\"26+52=78\", \"71-31=40\"],\n  [\"82-51=31\", \"13-13=0\"],\n  [\"15+41=56\", \"44+51=95\"],\n  [\"98-92=6\", \"72-49=23\"],\n  [\"38+12=50\", \"75-11=64\"],\n  [\"36-2=34\", \"15+59=74\"],\n  [\"56+42=98\", \"46+17=63\"],\n  [\"44+34=78\", \"52-41=11\"],\n  [\"77+7=84\", \"62-51=11\"],\n  [\"35+5=40\", \"63-51=12\"],\n  [\"86-0=86\", \"5+17=22\"],\n  [\"20+71=91\", \"96+3=99\"],\n  [\"55+15=70\", \"11+72=83\"],\n  [\"80+4=84\", \"72-55=17\"],\n  [\"31+41=72\", \"78-67=11\"],\n  [\"40+14=54\", \"78-22=56\"],\n  [\"16-14=2\", \"20-7=13\"],\n  [\"25-20=5\", \"53-47=6\"],\n  [\"56-30=26\", \"84-54=30\"],\n  [\"51-37=14\", \"67-26=41\"],\n  [\"10+27=37\", \"11+46=57\"],\n  [\"62-16=46\", \"54-28=26\"],\n  [\"10+81=91\", \"33+61=94\"],\n  [\"31+23=54\", \"50+25=75\"],\n  [\"29+3=32\", \"11+41=52\"],\n  [\"27+23=50\", \"60+24=84\"],\n  [\"49+27=76\", \"55+26=81\"],\n  [\"2+94=96\", \"54-10=44\"],\n  [\"17+72=89\", \"9+81=90\"],\n  [\"61-45=16\", \"26+69=95\"],\n  [\"31-24=7\", \"31+17=48\"],\n  [\"85-62=23\", \"82-2=80\"],\n  [\"32-10=22\", \"49-35=14\"],\n  [\"84-58=26\", \"37+35=72\"],\n  [\"64-63=1\", \"48+36=84\"],\n  [\"97-3=94\", \"99-21=78\"],\n  [\"26+21=47\", \"49-42=7\"],\n  [\"38-12=26\", \"41+50=91\"],\n  [\"99-88=11\", \"61-43=18\"],\n  [\"23+24=47\", \"53+4=57\"],\n  [\"32+22=54\", \"67+7=74\"],\n  [\"55+21=76\", \"4-2=2\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"53-23=30\", \"30-4=26\"),\n  @(\"76-35=41\", \"6+71=77\"),\n  @(\"35+55=90\", \"81-56=25\"),\n  @(\"20+38=58\", \"93-27=66\"),\n  @(\"40-16=24\", \"77-55=22\"),\n  @(\"39+8=47\", \"71-3=68\"),\n  @(\"73-33=40\", \"43-14=29\"),\n  @(\"2+56=58\", \"87+5=92\"),\n  @(\"4+92=96\", \"99-78=21\"),\n  @(\"30-29=1\", \"16+60=76\"),\n  @(\"95-56=39\", \"56-50=6\"),\n  @(\"76-21=55\", \"76-60=16\"),\n  @(\"74-12=62\", \"84-14=70\"),\n  @(\"57-10=47\", \"23+62=85\"),\n  @(\"50-15=35\", \"12+74=86\"),\n  @(\"22+29=51\", \"18+19=37\"),\n  @(\"95-12=83\", \"37+55=92\"),\n  @(\"70+11=81\", \"85-30=55\"),\n  @(\"24+34=58\", \"97-62=35\"),\n  @(\"41+53=94\", \"92-55=37\"),\n  @(\"73-1=72\", \"4+65=69\"),\n  @(\"63+33=96\", \"61-20=41\"),\n  @(\"43-13=30\", \"84-3=81\"),\n  @(\"13+31=44\", \"50-45=5\"),\n  @(\"57-22=35\", \"26+48=74\"),\n  @(\"67-30=37\", \"39-31=8\"),\n  @(\"1+42=43\", \"88-57=31\"),\n  @(\"12+72=84\", \"17+56=73\"),\n  @(\"87-82=5\", \"16+71=87\"),\n  @(\"75+11=86\", \"11+72=83\"),\n  @(\"7+53=60\", \"93-9=84\"),\n  @(\"60+35=95\", \"60-38=22\"),\n  @(\"64-51=13\", \"72-22=50\"),\n  @(\"66-44=22\", \"5+6=11\"),\n  @(\"46-32=14\", \"88-21=67\"),\n  @(\"20-3=17\", \"21+53=74\"),\n  @(\"12+86=98\", \"43+11=54\"),\n  @(\"4+25=29\", \"7+11=18\"),\n  @(\"38+23=61\", \"11+45=56\"),\n  @(\"33-16=17\", \"28-9=19\"),\n  @(\"69+1=70\", \"39-5=34\"),\n  @(\"28+25=53\", \"28+29=57\"),\n  @(\"58-55=3\", \"41-7=34\"),\n  @(\"22+11=33\", \"58-54=4\"),\n  @(\"21+70=91\", \"69-32=37\"),\n  @(\"31-25=6\", \"30+11=41\"),\n  @(\"66-18=48\", \"43+20=63\"),\n  @(\"53+28=81\", \"31+21=52\"),\n  @(\"5+58=63\", \"64-38=26\"),\n  @(\"18+38=56\", \"51-27=24\"),\n  @(\"60+26=86\", \"42+0=42\"),\n  @(\"19+45=64\", \"63-20=43\"),\n  @(\"64+14=78\", \"46-6=40\"),\n  @(\"26+30=56\", \"7+26=33\"),\n  @(\"47-30=17\", \"53-33=20\"),\n  @(\"68-38=30\", \"10+22=32\"),\n  @(\"36-20=16\", \"12+36=48\"),\n  @(\"29+59=88\", \"28+0=28\"),\n  @(\"26+52=78\", \"71-31=40\"),\n  @(\"82-51=31\", \"13-13=0\"),\n  @(\"15+41=56\", \"44+51=95\"),\n  @(\"98-92=6\", \"72-49=23\"),\n  @(\"38+12=50\", \"75-11=64\"),\n  @(\"36-2=34\", \"15+59=74\"),\n  @(\"56+42=98\", \"46+17=63\"),\n  @(\"44+34=78\", \"52-41=11\"),\n  @(\"77+7=84\", \"62-51=11\"),\n  @(\"35+5=40\", \"63-51=12\"),\n  @(\"86-0=86\", \"5+17=22\"),\n  @(\"20+71=91\", \"96+3=99\"),\n  @(\"55+15=70\", \"11+72=83\"),\n  @(\"80+4=84\", \"72-55=17\"),\n  @(\"31+41=72\", \"78-67=11\"),\n  @(\"40+14=54\", \"78-22=56\"),\n  @(\"16-14=2\", \"20-7=13\"),\n  @(\"25-20=5\", \"53-47=6\"),\n  @(\"56-30=26\", \"84-54=30\"),\n  @(\"51-37=14\", \"67-26=41\"),\n  @(\"10+27=37\", \"11+46=57\"),\n  @(\"62-16=46\", \"54-28=26\"),\n  @(\"10+81=91\", \"33+61=94\"),\n  @(\"31+23=54\", \"50+25=75\"),\n  @(\"29+3=32\", \"11+41=52\"),\n  @(\"27+23=50\", \"60+24=84\"),\n  @(\"49+27=76\", \"55+26=81\"),\n  @(\"2+94=96\", \"54-10=44\"),\n  @(\"17+72=89\", \"9+81=90\"),\n  @(\"61-45=16\", \"26+69=95\"),\n  @(\"31-24=7\", \"31+17=48\"),\n  @(\"85-62=23\", \"82-2=80\"),\n  @(\"32-10=22\", \"49-35=14\"),\n  @(\"84-58=26\", \"37+35=72\"),\n  @(\"64-63=1\", \"48+36=84\"),\n  @(\"97-3=94\", \"99-21=78\"),\n  @(\"26+21=47\", \"49-42=7\"),\n  @(\"38-12=26\", \"41+50=91\"),\n  @(\"99-88=11\", \"61-43=18\"),\n  @(\"23+24=47\", \"53+4=57\"),\n  @(\"32+22=54\", \"67+7=74\"),\n  @(\"55+21=76\", \"4-2=2\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}"}
